# Auto-generated script to update market-price-derived columns (H-N)
# across all 8 crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the upstream data refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5587.778
$ws.Range("I40").Value = 6026
$ws.Range("K40").Value = 6026
$ws.Range("M40").Value = -5851
$ws.Range("H98").Value = 1952.52
$ws.Range("I98").Value = 2058.5
$ws.Range("K98").Value = 2058.5
$ws.Range("M98").Value = -560.5
$ws.Range("H111").Value = 3111.6155
$ws.Range("I111").Value = 3223.375
$ws.Range("K111").Value = 9670.125
$ws.Range("M111").Value = -6603.125
$ws.Range("H122").Value = 1952.52
$ws.Range("I122").Value = 2058.5
$ws.Range("K122").Value = 6175.5
$ws.Range("M122").Value = -3725.5
$ws.Range("H131").Value = 7038.364
$ws.Range("I131").Value = 7131.6665
$ws.Range("K131").Value = 21394.9995
$ws.Range("M131").Value = -16354.9995
$ws.Range("H132").Value = 14433.82
$ws.Range("I132").Value = 1684.5676
$ws.Range("K132").Value = 5053.7028
$ws.Range("M132").Value = -2523.7028
$ws.Range("H135").Value = 11117135
$ws.Range("I135").Value = 16670023
$ws.Range("J135").Value = 11359.533
$ws.Range("K135").Value = 150030207
$ws.Range("L135").Value = 102235.797
$ws.Range("M135").Value = -150027672
$ws.Range("N135").Value = -107305.797

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3194.7576
$ws.Range("I2").Value = 2929.12
$ws.Range("J2").Value = 4024.875
$ws.Range("K2").Value = 2929.12
$ws.Range("L2").Value = 4024.875
$ws.Range("M2").Value = -2816.12
$ws.Range("N2").Value = -4250.875
$ws.Range("H45").Value = 4470
$ws.Range("J45").Value = 6726.857
$ws.Range("L45").Value = 6726.857
$ws.Range("N45").Value = -7480.857
$ws.Range("H61").Value = 2437.8462
$ws.Range("I61").Value = 2472.0454
$ws.Range("K61").Value = 2472.0454
$ws.Range("M61").Value = -2260.0454
$ws.Range("H74").Value = 3205.2273
$ws.Range("I74").Value = 2780.6
$ws.Range("J74").Value = 4115.143
$ws.Range("K74").Value = 2780.6
$ws.Range("L74").Value = 4115.143
$ws.Range("M74").Value = -1906.6
$ws.Range("N74").Value = -5863.143
$ws.Range("H77").Value = 3205.2273
$ws.Range("I77").Value = 2780.6
$ws.Range("J77").Value = 4115.143
$ws.Range("K77").Value = 13903
$ws.Range("L77").Value = 20575.715
$ws.Range("M77").Value = -9535
$ws.Range("N77").Value = -29311.715
$ws.Range("H116").Value = 3194.7576
$ws.Range("I116").Value = 2929.12
$ws.Range("J116").Value = 4024.875
$ws.Range("K116").Value = 2929.12
$ws.Range("L116").Value = 4024.875
$ws.Range("M116").Value = -635.1199999999999
$ws.Range("N116").Value = -8612.875
$ws.Range("H132").Value = 1366.2222
$ws.Range("I132").Value = 1388.3529
$ws.Range("J132").Value = 990
$ws.Range("K132").Value = 4165.0587
$ws.Range("L132").Value = 2970
$ws.Range("M132").Value = -1635.0587
$ws.Range("N132").Value = -8030
$ws.Range("H136").Value = 2437.8462
$ws.Range("I136").Value = 2472.0454
$ws.Range("K136").Value = 7416.1362
$ws.Range("M136").Value = -4866.1362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3194.7576
$ws.Range("I3").Value = 2929.12
$ws.Range("J3").Value = 4024.875
$ws.Range("K3").Value = 2929.12
$ws.Range("L3").Value = 4024.875
$ws.Range("M3").Value = -2815.12
$ws.Range("N3").Value = -4252.875
$ws.Range("H80").Value = 323.52942
$ws.Range("I80").Value = 240
$ws.Range("J80").Value = 334.66666
$ws.Range("K80").Value = 240
$ws.Range("L80").Value = 334.66666
$ws.Range("M80").Value = 758
$ws.Range("N80").Value = -2330.66666
$ws.Range("H83").Value = 323.52942
$ws.Range("I83").Value = 240
$ws.Range("J83").Value = 334.66666
$ws.Range("K83").Value = 1200
$ws.Range("L83").Value = 1673.3333
$ws.Range("M83").Value = 3792
$ws.Range("N83").Value = -11657.3333
$ws.Range("H134").Value = 8741.08
$ws.Range("I134").Value = 2496.4
$ws.Range("K134").Value = 7489.200000000001
$ws.Range("M134").Value = -4954.200000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6796.95
$ws.Range("J31").Value = 3081.6667
$ws.Range("L31").Value = 3081.6667
$ws.Range("N31").Value = -3671.6667
$ws.Range("H34").Value = 6796.95
$ws.Range("J34").Value = 3081.6667
$ws.Range("L34").Value = 3081.6667
$ws.Range("N34").Value = -3485.6667
$ws.Range("H88").Value = 33347.332
$ws.Range("J88").Value = 33347.332
$ws.Range("L88").Value = 33347.332
$ws.Range("N88").Value = -34159.332
$ws.Range("H91").Value = 33347.332
$ws.Range("J91").Value = 33347.332
$ws.Range("L91").Value = 33347.332
$ws.Range("N91").Value = -36155.332
$ws.Range("H99").Value = 8981803
$ws.Range("I99").Value = 1529747.6
$ws.Range("K99").Value = 1529747.6
$ws.Range("M99").Value = -1528249.6
$ws.Range("H123").Value = 29769.23
$ws.Range("J123").Value = 29769.23
$ws.Range("L123").Value = 29769.23
$ws.Range("N123").Value = -39569.23
$ws.Range("H126").Value = 8981803
$ws.Range("I126").Value = 1529747.6
$ws.Range("K126").Value = 4589242.800000001
$ws.Range("M126").Value = -4586772.800000001
$ws.Range("H132").Value = 7942.3687
$ws.Range("I132").Value = 8780.6875
$ws.Range("K132").Value = 26342.0625
$ws.Range("M132").Value = -23812.0625
$ws.Range("H134").Value = 2549.1924
$ws.Range("I134").Value = 2549.1924
$ws.Range("K134").Value = 7647.5772
$ws.Range("M134").Value = -5112.5772

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.944443
$ws.Range("I2").Value = 18.142857
$ws.Range("K2").Value = 108.857142
$ws.Range("M2").Value = 4.142858000000004
$ws.Range("H4").Value = 51350570
$ws.Range("I4").Value = 60412084
$ws.Range("K4").Value = 181236252
$ws.Range("M4").Value = -181236140
$ws.Range("H8").Value = 487.05554
$ws.Range("I8").Value = 487.05554
$ws.Range("K8").Value = 1461.16662
$ws.Range("M8").Value = -1322.16662
$ws.Range("H38").Value = 862.0968
$ws.Range("I38").Value = 63
$ws.Range("K38").Value = 189
$ws.Range("M38").Value = 158
$ws.Range("H106").Value = 7225
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 24000
$ws.Range("N106").Value = -25892
$ws.Range("H137").Value = 1497.5
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 1495
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 4485
$ws.Range("M137").Value = 600
$ws.Range("N137").Value = -14685

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H80").Value = 53183.363
$ws.Range("I80").Value = 141260
$ws.Range("J80").Value = 2853.8572
$ws.Range("K80").Value = 141260
$ws.Range("L80").Value = 2853.8572
$ws.Range("M80").Value = -140262
$ws.Range("N80").Value = -4849.8572
$ws.Range("H83").Value = 53183.363
$ws.Range("I83").Value = 141260
$ws.Range("J83").Value = 2853.8572
$ws.Range("K83").Value = 706300
$ws.Range("L83").Value = 14269.286
$ws.Range("M83").Value = -701308
$ws.Range("N83").Value = -24253.286
$ws.Range("H132").Value = 7254
$ws.Range("I132").Value = 7143.2085
$ws.Range("J132").Value = 7785.8
$ws.Range("K132").Value = 21429.6255
$ws.Range("L132").Value = 23357.4
$ws.Range("M132").Value = -18899.6255
$ws.Range("N132").Value = -28417.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 333.8
$ws.Range("I55").Value = 474.58334
$ws.Range("J55").Value = 122.625
$ws.Range("K55").Value = 474.58334
$ws.Range("L55").Value = 122.625
$ws.Range("M55").Value = -301.58334
$ws.Range("N55").Value = -468.625
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 53333.332
$ws.Range("J70").Value = 53333.332
$ws.Range("L70").Value = 53333.332
$ws.Range("N70").Value = -53963.332
$ws.Range("H73").Value = 53333.332
$ws.Range("J73").Value = 53333.332
$ws.Range("L73").Value = 53333.332
$ws.Range("N73").Value = -55517.332
$ws.Range("H104").Value = 55555.332
$ws.Range("J104").Value = 55555.332
$ws.Range("L104").Value = 55555.332
$ws.Range("N104").Value = -62543.332
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800
$ws.Range("H126").Value = 2734.1738
$ws.Range("I126").Value = 2618.3809
$ws.Range("K126").Value = 7855.1427
$ws.Range("M126").Value = -5385.1427
$ws.Range("H132").Value = 4026.432
$ws.Range("I132").Value = 3826
$ws.Range("K132").Value = 11478
$ws.Range("M132").Value = -8948
$ws.Range("H136").Value = 4082.2222
$ws.Range("I136").Value = 2617.84
$ws.Range("K136").Value = 7853.52
$ws.Range("M136").Value = -5303.52
